$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.823570013046265
$ws.Range("B1").Value = 5.86894416809082
$ws.Range("C1").Value = 4.935348987579346
$ws.Range("D1").Value = 5.746235847473145
$ws.Range("E1").Value = 4.157257080078125
